$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the organization name value from B1 (cell is cleared/removed entirely)
$ws.Range("B1").ClearContents()

# Update the "Objetivo de Gobierno" text in B8
$ws.Range("B8").Value = "Generar condiciones de paz y tranquilidad para la ciudadanía, preservando en todo momento el Estado de derecho y cumpliendo el mandato de brindar protección a las personas y sus bienes, garantizando el derecho de acceso a la justicia, promoviendo la mediación comunitaria y aplicando la ley a cabalidad, privilegiando la rendición de cuentas y la participación comunitaria."

# Update the report generation date/time in B12
$ws.Range("B12").Value = "07-07-2022 12:22:48 pm"

# Helper to write a numeric-looking value as plain text (matching the
# source file, which stores these as inline strings, not numbers)
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Add the new "Fin" indicator row (row 17)
$ws.Range("A17").Value = "Fin"
$ws.Range("B17").Value = "2.1 Aumento en el nivel educativo promedio de la población del estado. "
$ws.Range("C17").Value = "gestion"
$ws.Range("D17").Value = "Eficacia"
$ws.Range("E17").Value = "Agenda de Gobierno Digital"
Set-TextValue "F17" "5512"
$ws.Range("G17").Value = "Porcentaje de Proyectos del Plan Querétaro Digital liberados con respecto al total de Proyectos"
Set-TextValue "H17" "0"
$ws.Range("I17").Value = "Trimestral"
$ws.Range("J17").Value = "Sumatoria"
$ws.Range("K17").Value = "A"
$ws.Range("L17").Value = "Proyecto"
$ws.Range("M17").Value = "A"
$ws.Range("N17").Value = "wefewf"
Set-TextValue "O17" "10"
Set-TextValue "P17" "0"
Set-TextValue "Q17" "0"
$ws.Range("AA17").Value = 10
